$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1379.6666
$ws.Range("I12").Value = 140
$ws.Range("J12").Value = 1999.5
$ws.Range("K12").Value = 140
$ws.Range("L12").Value = 1999.5
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = -2339.5

$ws.Range("H51").Value = 11804.286
$ws.Range("J51").Value = 7105
$ws.Range("L51").Value = 7105
$ws.Range("N51").Value = -8073

$ws.Range("H61").Value = 2522.2222
$ws.Range("I61").Value = 2600
$ws.Range("K61").Value = 7800
$ws.Range("M61").Value = -7628

$ws.Range("H80").Value = 683.1818
$ws.Range("I80").Value = 396.33334
$ws.Range("J80").Value = 1027.4
$ws.Range("K80").Value = 1189.00002
$ws.Range("L80").Value = 3082.2
$ws.Range("M80").Value = -191.0000199999999
$ws.Range("N80").Value = -5078.200000000001

$ws.Range("H83").Value = 683.1818
$ws.Range("I83").Value = 396.33334
$ws.Range("J83").Value = 1027.4
$ws.Range("K83").Value = 3567.00006
$ws.Range("L83").Value = 9246.6
$ws.Range("M83").Value = 1424.99994
$ws.Range("N83").Value = -19230.6

$ws.Range("H87").Value = 30000
$ws.Range("I87").Value = 30000
$ws.Range("K87").Value = 30000
$ws.Range("M87").Value = -28752

$ws.Range("H90").Value = 30000
$ws.Range("I90").Value = 30000
$ws.Range("K90").Value = 90000
$ws.Range("M90").Value = -83760

$ws.Range("H92").Value = 4153.5713
$ws.Range("I92").Value = 3895
$ws.Range("K92").Value = 3895
$ws.Range("M92").Value = -2647

$ws.Range("H113").Value = 7098.7144
$ws.Range("J113").Value = 7815.25
$ws.Range("L113").Value = 7815.25
$ws.Range("N113").Value = -14323.25

$ws.Range("H131").Value = 5760.25
$ws.Range("I131").Value = 2713.6667
$ws.Range("K131").Value = 8141.000100000001
$ws.Range("M131").Value = -3101.000100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2862.5908
$ws.Range("I61").Value = 2237.7222
$ws.Range("J61").Value = 5674.5
$ws.Range("K61").Value = 2237.7222
$ws.Range("L61").Value = 5674.5
$ws.Range("M61").Value = -2025.7222
$ws.Range("N61").Value = -6098.5

$ws.Range("H74").Value = 1467.2354
$ws.Range("I74").Value = 635.0345
$ws.Range("K74").Value = 635.0345
$ws.Range("M74").Value = 238.9655

$ws.Range("H77").Value = 1467.2354
$ws.Range("I77").Value = 635.0345
$ws.Range("K77").Value = 3175.1725
$ws.Range("M77").Value = 1192.8275

$ws.Range("H136").Value = 2862.5908
$ws.Range("I136").Value = 2237.7222
$ws.Range("J136").Value = 5674.5
$ws.Range("K136").Value = 6713.1666
$ws.Range("L136").Value = 17023.5
$ws.Range("M136").Value = -4163.1666
$ws.Range("N136").Value = -22123.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 955.5714
$ws.Range("I20").Value = 938.2
$ws.Range("J20").Value = 999
$ws.Range("K20").Value = 938.2
$ws.Range("L20").Value = 999
$ws.Range("M20").Value = -691.2
$ws.Range("N20").Value = -1493

$ws.Range("H107").Value = 2210.5454
$ws.Range("I107").Value = 1536.5883
$ws.Range("K107").Value = 1536.5883
$ws.Range("M107").Value = 383.4117000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H86").Value = 3630.652
$ws.Range("I86").Value = 3613.1538
$ws.Range("J86").Value = 3653.4
$ws.Range("K86").Value = 3613.1538
$ws.Range("L86").Value = 3653.4
$ws.Range("M86").Value = -2490.1538
$ws.Range("N86").Value = -5899.4

$ws.Range("H89").Value = 3630.652
$ws.Range("I89").Value = 3613.1538
$ws.Range("J89").Value = 3653.4
$ws.Range("K89").Value = 18065.769
$ws.Range("L89").Value = 18267
$ws.Range("M89").Value = -12449.769
$ws.Range("N89").Value = -29499

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").ClearContents()

$ws.Range("H105").Value = 1194.7307
$ws.Range("I105").Value = 1155
$ws.Range("K105").Value = 1155
$ws.Range("M105").Value = 592

$ws.Range("H107").Value = 974
$ws.Range("J107").Value = 949
$ws.Range("L107").Value = 949
$ws.Range("N107").Value = -4789

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5651854.5
$ws.Range("I11").Value = 5651854.5
$ws.Range("K11").Value = 16955563.5
$ws.Range("M11").Value = -16955423.5

$ws.Range("H87").Value = 13002.8
$ws.Range("I87").Value = 11671.333
$ws.Range("K87").Value = 35013.999
$ws.Range("M87").Value = -33765.999

$ws.Range("H90").Value = 13002.8
$ws.Range("I90").Value = 11671.333
$ws.Range("K90").Value = 105041.997
$ws.Range("M90").Value = -98801.997

$ws.Range("H99").Value = 26000
$ws.Range("I99").Value = 22500
$ws.Range("K99").Value = 67500
$ws.Range("M99").Value = -65254

$ws.Range("H113").Value = 651
$ws.Range("J113").Value = 439.5
$ws.Range("L113").Value = 1318.5
$ws.Range("N113").Value = -5658.5

$ws.Range("H131").Value = 1963.2858
$ws.Range("J131").Value = 2046
$ws.Range("L131").Value = 6138
$ws.Range("N131").Value = -16218

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3471.111
$ws.Range("I80").Value = 3397
$ws.Range("K80").Value = 3397
$ws.Range("M80").Value = -2399

$ws.Range("H83").Value = 3471.111
$ws.Range("I83").Value = 3397
$ws.Range("K83").Value = 16985
$ws.Range("M83").Value = -11993

$ws.Range("H122").Value = 3515.9
$ws.Range("I122").Value = 2907.3
$ws.Range("J122").Value = 4733.1
$ws.Range("K122").Value = 8721.900000000001
$ws.Range("L122").Value = 14199.3
$ws.Range("M122").Value = -6271.900000000001
$ws.Range("N122").Value = -19099.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1890.8334
$ws.Range("J22").Value = 2499.3333
$ws.Range("L22").Value = 2499.3333
$ws.Range("N22").Value = -3089.3333

$ws.Range("H27").Value = 1890.8334
$ws.Range("J27").Value = 2499.3333
$ws.Range("L27").Value = 2499.3333
$ws.Range("N27").Value = -2713.3333

$ws.Range("H40").Value = 2354.389
$ws.Range("I40").Value = 2115.3333
$ws.Range("J40").Value = 2832.5
$ws.Range("K40").Value = 2115.3333
$ws.Range("L40").Value = 2832.5
$ws.Range("M40").Value = -1979.3333
$ws.Range("N40").Value = -3104.5

$ws.Range("H46").Value = 1571.1666
$ws.Range("I46").Value = 1074.25
$ws.Range("K46").Value = 1074.25
$ws.Range("M46").Value = -886.25

$ws.Range("H55").Value = 430.66666
$ws.Range("I55").Value = 419.75
$ws.Range("J55").Value = 436.125
$ws.Range("K55").Value = 419.75
$ws.Range("L55").Value = 436.125
$ws.Range("M55").Value = -246.75
$ws.Range("N55").Value = -782.125

$ws.Range("H108").Value = 58966
$ws.Range("J108").Value = 58966
$ws.Range("L108").Value = 58966
$ws.Range("N108").Value = -66646

$ws.Range("H109").Value = 63944
$ws.Range("J109").Value = 63944
$ws.Range("L109").Value = 63944
$ws.Range("N109").Value = -66718

$ws.Range("H136").Value = 2145.2964
$ws.Range("I136").Value = 1780.7693
$ws.Range("J136").Value = 2483.7856
$ws.Range("K136").Value = 5342.3079
$ws.Range("L136").Value = 7451.3568
$ws.Range("M136").Value = -2792.3079
$ws.Range("N136").Value = -12551.3568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 69377
$ws.Range("J109").Value = 69377
$ws.Range("L109").Value = 69377
$ws.Range("N109").Value = -72151

$ws.Range("H117").Value = 44999
$ws.Range("J117").Value = 44999
$ws.Range("L117").Value = 44999
$ws.Range("N117").Value = -54177
